# Commit: "add consump to stock"
#
# The underlying change consolidates the two near-duplicate row labels
# "لوله های گازی و صنعتی" and "انواع لوله گازی و صنعتی" (both referring to
# gas/industrial pipe product lines) into a single shorter label "لوله"
# ("pipe") everywhere they occur in the "Overview" sheet, and also switches
# the sheet to a right-to-left view with B17 as the active/selected cell
# (matching how the author was last working in the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Cells whose label needs to become "لوله" (pipe) instead of the old,
# longer "لوله های گازی و صنعتی" / "انواع لوله گازی و صنعتی" labels.
$targetCells = @("B13", "B17", "B18", "B34", "B51", "B68", "B84", "B100")

foreach ($addr in $targetCells) {
    $ws.Range($addr).Value = "لوله"
}

# Switch the sheet view to right-to-left (matches the Farsi content) and
# restore the author's last selection (B17).
$ws.Activate()
$excel.ActiveWindow.DisplayRightToLeft = $true
$ws.Range("B17").Select()
